$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before column B ("SPECIES_CN"), shifting old B..O to C..P.
$ws.Columns("B").Insert()

# 2. Insert a new row before row 5 (the new "绿樟" record), shifting old row 5 (now row 6)
#    and below down by one.
$ws.Rows(5).Insert()

# 3. Relabel the header row.
$ws.Range("A1").Value = "YOUR_SEARCH"
$ws.Range("B1").Value = "SPECIES_CN"

# 4. Populate the new SPECIES_CN column by duplicating the (TAXA_NAME/YOUR_SEARCH) column
#    for each existing data row.
$ws.Range("B2").Value = $ws.Range("A2").Text
$ws.Range("B3").Value = $ws.Range("A3").Text
$ws.Range("B4").Value = $ws.Range("A4").Text
$ws.Range("B6").Value = $ws.Range("A6").Text

# 5. Fill in the newly inserted row with the single known value; the rest of the
#    row is left blank (materialised as empty cells, matching the inserted row).
$ws.Range("A5").Value = "绿樟"
$ws.Range("B5:P5").Style = "Normal"
